# Insert a new data row for Jengibre (Vega Modelo de Temuco) at row 28.
# This pushes all existing data rows (28..134) down by one (becoming 29..135).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(28).Insert()

$ws.Cells.Item(28, 1).Value2 = 10
$ws.Cells.Item(28, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(28, 3).Value = "La Araucanía"
$ws.Cells.Item(28, 4).Value2 = 44620
$ws.Cells.Item(28, 5).Value2 = 9
$ws.Cells.Item(28, 6).Value2 = 100114007
$ws.Cells.Item(28, 7).Value = "Jengibre"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value2 = 20
$ws.Cells.Item(28, 11).Value2 = 26000
$ws.Cells.Item(28, 12).Value2 = 26000
$ws.Cells.Item(28, 13).Value2 = 26000
$ws.Cells.Item(28, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(28, 15).Value = "Perú"
$ws.Cells.Item(28, 16).Value2 = 2000
$ws.Cells.Item(28, 17).Value2 = 13
$ws.Cells.Item(28, 18).Value = "Hortaliza"
